$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clear cells that are being vacated (old "Tests for Enum" block at rows
#    8-12 is being replaced by new content at rows 6-9 and moved down to
#    rows 11-15).
# ---------------------------------------------------------------------------
$ws.Range("A10:E10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("B11:E11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("G12:H12").ClearContents()
$ws.Range("J12:K12").ClearContents()

# ---------------------------------------------------------------------------
# 2) New shared strings must be appended to the shared-string table in a
#    specific order (matching the target workbook) - write those cells
#    first, in this exact sequence.
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "empty cells in input array are considered as 0"
$ws.Range("M7").Value = "a single input cell with #NA! is considered as null"
$ws.Range("M8").Value = "but an input array full of #NA! should cause a conversion error"
$ws.Range("J3").Value = "Inputs"

# ---------------------------------------------------------------------------
# 3) Row 1 - B1 now references the extended range A4:A15 and evaluates TRUE.
# ---------------------------------------------------------------------------
$ws.Range("B1").Formula = "=AND(A4:A15)"

# ---------------------------------------------------------------------------
# 4) Row 3 header row - new "Inputs" header cell, bold like its neighbours.
# ---------------------------------------------------------------------------
$ws.Range("J3").Style = $ws.Range("A3").Style

# ---------------------------------------------------------------------------
# 5) "Tests for Complex" block (rows 4-9)
# ---------------------------------------------------------------------------

# Row 4 unchanged in content, but B4:C9 become one shared formula block.
$ws.Range("B4:C9").Formula = "=D4-G4"
$ws.Range("C4").Formula = "=E4-H4"
$ws.Range("C9").Formula = "=E9-H9"

# Row 5 gets a new label in M5 (reuse of existing shared string).
$ws.Range("M5").Value = "null value provided as missing parameters"

# Row 6: dnaNullableComplex now takes J6:K6 (currently-empty cells) instead
# of G6:H6, so the previously-error results become 0, and G6/H6 become 0.
$ws.Range("D6:E6").FormulaArray = "=_xll.dnaNullableComplex(J6:K6)"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# Row 7 (new row) - single-cell #N/A input is treated as null.
$ws.Range("A7").Formula = "=SUMSQ(B7:C7)<0.00000000000001"
$ws.Range("B7").Formula = "=D7-G7"
$ws.Range("C7").Formula = "=E7-H7"
$ws.Range("D7:E7").FormulaArray = "=_xll.dnaNullableComplex(J7)"
$ws.Range("E7").Value = 222
$ws.Range("G7").Value = 111
$ws.Range("H7").Value = 222
$ws.Range("J7").Formula = "=NA()"

# Row 8 (new row) - input array full of #N/A causes a conversion error.
$ws.Range("A8").Formula = "=SUMSQ(B8:C8)<0.00000000000001"
$ws.Range("B8").Formula = "=D8-G8"
$ws.Range("C8").Formula = "=E8-H8"
$ws.Range("D8:E8").FormulaArray = "=ERROR.TYPE(_xll.dnaNullableComplex(J8:K8))"
$ws.Range("E8").Value = 6
$ws.Range("G8").Formula = "=ERROR.TYPE(#NUM!)"
$ws.Range("H8").Value = 6
$ws.Range("J8").Formula = "=NA()"
$ws.Range("K8").Formula = "=NA()"

# Row 9: now exercises dnaNullableComplex(G9:H9) instead of dnaNullableEnum().
$ws.Range("D9:E9").FormulaArray = "=_xll.dnaNullableComplex(G9:H9)"
$ws.Range("E9").Value = 4
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 4
$ws.Range("M9").Value = "non-null values"

# ---------------------------------------------------------------------------
# 6) "Tests for Enum" block, now at rows 11-15.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Tests for Enum"
$ws.Range("A11").Style = $ws.Range("A3").Style
$ws.Range("G11").Value = "Expected values"
$ws.Range("G11").Style = $ws.Range("A3").Style
$ws.Range("J11").Value = "Inputs"
$ws.Range("J11").Style = $ws.Range("A3").Style

$ws.Range("A12:A15").Formula = "=SUMSQ(B12:C12)<0.00000000000001"
$ws.Range("B12:C15").Formula = "=D12-G12"
$ws.Range("C12").Formula = "=E12-H12"
$ws.Range("C13").Formula = "=E13-H13"
$ws.Range("C14").Formula = "=E14-H14"
$ws.Range("C15").Formula = "=E15-H15"

$ws.Range("D12:E12").FormulaArray = "=_xll.dnaNullableEnum()"
$ws.Range("E12").Value = 0
$ws.Range("M12").Value = "null value provided as missing parameters"

$ws.Range("D13:E13").FormulaArray = "=_xll.dnaNullableEnum(J13,K13)"
$ws.Range("E13").Value = 0
$ws.Range("M13").Value = "null values provided as refs to empty cells"

$ws.Range("D14:E14").FormulaArray = "=_xll.dnaNullableEnum(J14,K14)"
$ws.Range("E14").Value = 0
$ws.Range("J14").Formula = "=NA()"
$ws.Range("K14").Formula = "=NA()"
$ws.Range("M14").Value = "null values provided as #N/As"

$ws.Range("D15:E15").FormulaArray = "=_xll.dnaNullableEnum(J15,K15)"
$ws.Range("E15").Value = 1
$ws.Range("G15").Value = -1
$ws.Range("H15").Value = 1
$ws.Range("J15").Value = "negative"
$ws.Range("K15").Value = "imaginary"
$ws.Range("M15").Value = "non-null values"

# ---------------------------------------------------------------------------
# 7) Selection moves to J11, matching the saved UI state.
# ---------------------------------------------------------------------------
$ws.Range("J11").Select()
